# Update the "想去人数" (column F) counts on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Map of worksheet name -> hashtable of row number -> new F value
$updates = @{
    "展览"   = @{
        5  = 1741
        6  = 3311
        8  = 2164
        11 = 589
        13 = 1652
        15 = 75
        18 = 182
        19 = 1545
        20 = 592
        21 = 696
        23 = 12124
        24 = 12136
        25 = 895
        30 = 323
        33 = 546
    }
    "全部类型" = @{
        6  = 1741
        7  = 3311
        9  = 2164
        12 = 589
        14 = 1652
        17 = 75
        22 = 182
        23 = 1545
        24 = 592
        25 = 696
        27 = 12124
        28 = 12136
        29 = 895
        34 = 323
        39 = 546
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($r in $rows.Keys) {
        $ws.Cells.Item($r, 6).Value = $rows[$r]
    }
}
